$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("A3").Value = "Otro"

$ws.Range("A4").Select()
